# feat: add 2022-Q4 data
#
# - Inserts a new "2022-Q4" sheet (copy of the previous quarter's sheet
#   layout) right after "总计", pushing 2022-Q3/2022-Q2/2022-Q1/2021-Q3
#   one position to the right.
# - Updates the "总计" summary sheet with a new top data row for 2022-Q4
#   and shifts the previous rows down by one.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item(2)

# --- 1. Create the new "2022-Q4" sheet as a copy of the old "2022-Q3"
#        sheet (keeps header row / cell types / styles identical), placed
#        right after "总计". ---
$q3Sheet.Copy($null, $summary) | Out-Null
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Update the per-fund figures for the new quarter. Columns D-G are stored
# as text in the source file (e.g. "3.80"), so the NumberFormat dance below
# stops Excel's autodetect from silently converting them to numbers; the
# PasteSpecial afterwards restores the (unstyled) look of the other cells
# in the row instead of leaving the temporary "@" text format behind.
function Set-TextValue($cell, $text, $formatDonor) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $formatDonor.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

Set-TextValue $q4Sheet.Range("D2") "3.80" $q4Sheet.Range("C2")
Set-TextValue $q4Sheet.Range("E2") "90.62" $q4Sheet.Range("C2")
Set-TextValue $q4Sheet.Range("F2") "3.68" $q4Sheet.Range("C2")
Set-TextValue $q4Sheet.Range("G2") "0.1398" $q4Sheet.Range("C2")
$q4Sheet.Range("H2").Value = 10

# --- 2. Update the "总计" summary sheet: insert a new first data row for
#        2022-Q4 and push the rest down by one. ---
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.14

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.17

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 3
$summary.Range("D4").Value = 0.22

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0.03

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 2
$summary.Range("D6").Value = 0.03

# Copy the "A" column header style down onto the newly-used A6 cell so it
# matches the bordered/bold look of the rest of column A.
$summary.Range("A5").Copy() | Out-Null
$summary.Range("A6").PasteSpecial(-4122) | Out-Null
$summary.Range("A6").Value = 4

# --- 3. Restore the originally-selected tab (2021-Q3, the last sheet). ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Activate()

$excel.CutCopyMode = $false
